$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 112, shifting existing rows 112:149 down to 113:150
$ws.Rows.Item(112).Insert()

# Fill in the new row 112 with its data (columns constant across the dataset
# are copied from the surrounding rows; the varying columns take the new values)
$ws.Range("A112").Value = 9
$ws.Range("B112").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C112").Value = "Metropolitana"
$ws.Range("D112").Value = 44988
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100101
$ws.Range("H112").Value = "Berries"
$ws.Range("I112").Value = 100101004
$ws.Range("J112").Value = "Frambuesa"
$ws.Range("K112").Value = "Sin especificar"
$ws.Range("L112").Value = "Primera"
$ws.Range("M112").Value = 290
$ws.Range("N112").Value = 7000
$ws.Range("O112").Value = 7000
$ws.Range("P112").Value = 7000
$ws.Range("Q112").Value = "$/bandeja 2 kilos"
$ws.Range("R112").Value = "Región Metropolitana"
$ws.Range("S112").Value = 3500
$ws.Range("T112").Value = 2
